$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 12 ("Exercise: Library") - Content Placeholder 2
# ---------------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$shape12 = $s12.Shapes.Item(2)
$tr12 = $shape12.TextFrame.TextRange

# Paragraph 1: "Create a library application that is able to:"
# Run 3 ("that is able ", chars 30-42) needs to become two runs:
#   "where a user using a console interface " + "is able "
# Splitting is achieved by re-writing only the leading part of the run
# (chars 30-34, i.e. "that "), which forces the engine to split the run in
# two while leaving the remainder ("is able ") as its own run.
$para1Run3Lead = $tr12.Characters(30, 5)
$para1Run3Lead.Text = "where a user using a console interface "

# Paragraph 6 (last paragraph): "Use JDBC to store the data in " + "the <book> table..."
# merge the two runs of that paragraph into a single run.
$para6 = $tr12.Paragraphs(6, 1)
$para6.Text = "TEMP_MERGE_PLACEHOLDER"
$tr12b = $shape12.TextFrame.TextRange
$para6b = $tr12b.Paragraphs(6, 1)
$para6b.Text = "Use JDBC to store the data in the " + [char]0x201C + "book" + [char]0x201D + " table created in the previous exercise."

# ---------------------------------------------------------------------------
# Slide 13 ("Application: AliExpress") - Content Placeholder 2
# ---------------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$shape13 = $s13.Shapes.Item(2)
$tr13 = $shape13.TextFrame.TextRange

# Paragraph 1: merge 7 runs ("Provide a second implementation..." ... "memory.")
# into a single run.
$para13_1 = $tr13.Paragraphs(1, 1)
$para13_1.Text = "TEMP_MERGE_PLACEHOLDER"
$tr13b = $shape13.TextFrame.TextRange
$para13_1b = $tr13b.Paragraphs(1, 1)
$para13_1b.Text = "Provide a second implementation for the data access component(s) of the warehouse module that uses JDBC to store the products in a relational database rather than in memory."

# Paragraph 2: merge 3 runs ("Control which implementation..." ... "class.")
# into a single run, keeping the endParaRPr.
$tr13c = $shape13.TextFrame.TextRange
$para13_2 = $tr13c.Paragraphs(2, 1)
$para13_2.Text = "TEMP_MERGE_PLACEHOLDER"
$tr13d = $shape13.TextFrame.TextRange
$para13_2b = $tr13d.Paragraphs(2, 1)
$para13_2b.Text = "Control which implementation (in-memory or relational database) the application uses in the " + [char]0x201C + "main" + [char]0x201D + " application class."
